$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on cells whose new values look numeric, so Excel
# stores them as text (matching the original inlineStr cell type) instead
# of silently converting to a number and losing formatting (trailing zeros, etc).
$textCells = @("D5","D6","D9","D10","D11","D12","D13","D14","D19","D20","D21","D22","D23","D24","D25","D26","D27","D29","D30","D31","D32","D33","D36","D38","D39","D41","D43","D44","D45","D46","D47","D48","D49","D50")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated values
$ws.Range("D2").Value = "70.007.46"
$ws.Range("E2").Value = "  +3.94%  "
$ws.Range("D3").Value = "3.781.14"
$ws.Range("E3").Value = "  +21.41%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "615.39"
$ws.Range("E5").Value = "  +6.79%  "
$ws.Range("D6").Value = "176.37"
$ws.Range("E6").Value = "  -1.17%  "
$ws.Range("D7").Value = "3.777.05"
$ws.Range("E7").Value = "  +21.34%  "
$ws.Range("E8").Value = "  +0.09%  "
$ws.Range("D9").Value = "0.547"
$ws.Range("E9").Value = "  +5.90%  "
$ws.Range("D10").Value = "0.167"
$ws.Range("E10").Value = "  +9.73%  "
$ws.Range("D11").Value = "6.37"
$ws.Range("E11").Value = "  -2.54%  "
$ws.Range("D12").Value = "0.500"
$ws.Range("E12").Value = "  +7.10%  "
$ws.Range("D13").Value = "40.43"
$ws.Range("E13").Value = "  +10.83%  "
$ws.Range("D14").Value = "0.0000256"
$ws.Range("E14").Value = "  +6.19%  "
$ws.Range("D15").Value = "4.422.56"
$ws.Range("E15").Value = "  +21.75%  "
$ws.Range("D16").Value = "3.787.37"
$ws.Range("E16").Value = "  +21.72%  "
$ws.Range("D17").Value = "70.274.22"
$ws.Range("E17").Value = "  +4.39%  "
$ws.Range("E18").Value = "  +1.22%  "
$ws.Range("D19").Value = "7.59"
$ws.Range("E19").Value = "  +7.76%  "
$ws.Range("D20").Value = "524.40"
$ws.Range("E20").Value = "  +7.96%  "
$ws.Range("D21").Value = "16.68"
$ws.Range("E21").Value = "  +1.42%  "
$ws.Range("D22").Value = "9.44"
$ws.Range("E22").Value = "  +22.51%  "
$ws.Range("D23").Value = "0.744"
$ws.Range("E23").Value = "  +7.85%  "
$ws.Range("D24").Value = "88.55"
$ws.Range("E24").Value = "  +5.86%  "
$ws.Range("D25").Value = "2.48"
$ws.Range("E25").Value = "  +8.57%  "
$ws.Range("D26").Value = "13.50"
$ws.Range("E26").Value = "  +6.21%  "
$ws.Range("D27").Value = "10.90"
$ws.Range("E27").Value = "  +3.51%  "
$ws.Range("E28").Value = "  +0.01%  "
$ws.Range("D29").Value = "0.0000124"
$ws.Range("E29").Value = "  +31.19%  "
$ws.Range("D30").Value = "2.51"
$ws.Range("E30").Value = "  +8.11%  "
$ws.Range("D31").Value = "2.87"
$ws.Range("E31").Value = "  +9.55%  "
$ws.Range("D32").Value = "7.90"
$ws.Range("E32").Value = "  -1.15%  "
$ws.Range("D33").Value = "31.96"
$ws.Range("E33").Value = "  +13.83%  "
$ws.Range("E34").Value = "  +2.19%  "
$ws.Range("E35").Value = "  +0.10%  "
$ws.Range("D36").Value = "6.19"
$ws.Range("E36").Value = "  +10.69%  "
$ws.Range("E37").Value = "  +9.83%  "
$ws.Range("D38").Value = "0.342"
$ws.Range("E38").Value = "  +7.16%  "
$ws.Range("D39").Value = "0.133"
$ws.Range("E39").Value = "  +8.01%  "
$ws.Range("E40").Value = "  +7.06%  "
$ws.Range("D41").Value = "51.57"
$ws.Range("E41").Value = "  +4.74%  "
$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").Value = "3.137.68"
$ws.Range("E42").Value = "  +12.54%  "
$ws.Range("B43").Value = "Cosmos"
$ws.Range("C43").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D43").Value = "8.84"
$ws.Range("E43").Value = "  +6.40%  "
$ws.Range("B44").Value = "Arweave"
$ws.Range("C44").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D44").Value = "44.37"
$ws.Range("E44").Value = "  -7.82%  "
$ws.Range("B45").Value = "Bittensor"
$ws.Range("C45").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D45").Value = "426.80"
$ws.Range("E45").Value = "  +14.10%  "
$ws.Range("D46").Value = "2.73"
$ws.Range("E46").Value = "  +1.49%  "
$ws.Range("D47").Value = "0.0367"
$ws.Range("E47").Value = "  +5.98%  "
$ws.Range("D48").Value = "27.69"
$ws.Range("E48").Value = "  +3.03%  "
$ws.Range("B49").Value = "Monero"
$ws.Range("C49").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D49").Value = "137.44"
$ws.Range("E49").Value = "  +1.40%  "
$ws.Range("B50").Value = "ThetaToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D50").Value = "2.52"
$ws.Range("E50").Value = "  +6.03%  "

# Restore default (Normal) style on the cells we temporarily reformatted,
# now that the text values are committed, so no stray formatting remains.
foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
